# Update the "Comuna" column (D) on sheet INCO: replace the
# "<code> - <neighborhood name>" text with just the comuna number,
# keeping the cell as text (not a numeric value).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("INCO")

$comunaByRow = @{
    2  = "1"
    3  = "2"
    4  = "12"
    5  = "13"
    6  = "4"
    7  = "2"
    8  = "12"
    9  = "1"
    10 = "3"
    11 = "3"
    12 = "3"
    13 = "3"
    14 = "13"
    15 = "2"
    16 = "2"
    17 = "11"
    18 = "13"
    19 = "12"
    20 = "12"
    21 = "12"
    22 = "12"
    23 = "12"
    24 = "12"
    25 = "13"
    26 = "11"
    27 = "3"
    28 = "15"
    29 = "5"
    30 = "10"
    31 = "11"
    32 = "15"
    33 = "13"
    34 = "11"
    35 = "3"
}

foreach ($row in $comunaByRow.Keys) {
    $cell = $ws.Range("D$row")
    # Force text storage so "1", "2", "12", ... are not reinterpreted as numbers.
    $cell.NumberFormat = "@"
    $cell.Value = $comunaByRow[$row]
}
